$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "70.861.04"
$ws.Range("E2").Value = "  +1.97%  "

$ws.Range("D3").Value = "3.577.32"
$ws.Range("E3").Value = "  +2.19%  "

$ws.Range("E4").Value = "  +0.15%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "586.24"
$ws.Range("E5").Value = "  -0.18%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "189.21"
$ws.Range("E6").Value = "  +2.64%  "

$ws.Range("D7").Value = "3.571.18"
$ws.Range("E7").Value = "  +2.42%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.624"
$ws.Range("E8").Value = "  +2.29%  "

$ws.Range("E9").Value = "  -0.09%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.219"
$ws.Range("E10").Value = "  +10.28%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.649"
$ws.Range("E11").Value = "  +0.79%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "54.25"
$ws.Range("E12").Value = "  +1.24%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000313"
$ws.Range("E13").Value = "  +2.73%  "

$ws.Range("E14").Value = "  +0.92%  "

$ws.Range("D15").Value = "4.147.96"
$ws.Range("E15").Value = "  +2.22%  "

$ws.Range("D16").Value = "70.890.59"
$ws.Range("E16").Value = "  +2.09%  "

$ws.Range("D17").Value = "3.598.86"
$ws.Range("E17").Value = "  +2.82%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "19.12"
$ws.Range("E18").Value = "  -0.37%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.77"
$ws.Range("E19").Value = "  +4.63%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "573.55"
$ws.Range("E20").Value = "  +6.35%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.121"
$ws.Range("E21").Value = "  +0.89%  "

$ws.Range("E22").Value = "  -0.20%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "18.08"
$ws.Range("E23").Value = "  -2.19%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "4.61"
$ws.Range("E24").Value = "  +1.30%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "4.95"
$ws.Range("E25").Value = "  +2.06%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "94.40"
$ws.Range("E26").Value = "  -0.29%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "11.12"
$ws.Range("E27").Value = "  +0.68%  "

$ws.Range("E28").Value = "  -0.82%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.35"
$ws.Range("E29").Value = "  +3.42%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "32.63"
$ws.Range("E30").Value = "  +2.41%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.12"
$ws.Range("E31").Value = "  -2.12%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "12.25"
$ws.Range("E32").Value = "  -1.95%  "

$ws.Range("E33").Value = "  +3.08%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.90"
$ws.Range("E34").Value = "  +24.50%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "63.28"
$ws.Range("E35").Value = "  -1.17%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "3.25"
$ws.Range("E36").Value = "  +7.25%  "

$ws.Range("B37").Value = "Bittensor"
$ws.Range("C37").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "536.76"
$ws.Range("E37").Value = "  -0.74%  "

$ws.Range("B38").Value = "TheGraph"
$ws.Range("C38").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.412"
$ws.Range("E38").Value = "  +1.47%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "38.44"
$ws.Range("E39").Value = "  +2.44%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.00"
$ws.Range("E40").Value = "  +0.01%  "

$ws.Range("B41").Value = "Maker"
$ws.Range("C41").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D41").Value = "3.634.04"
$ws.Range("E41").Value = "  +9.19%  "

$ws.Range("B42").Value = "PEPE"
$ws.Range("C42").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D42").Value = "0.0₃0799"
$ws.Range("E42").Value = "  +4.91%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.140"
$ws.Range("E43").Value = "  +5.33%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "3.52"
$ws.Range("E44").Value = "  +4.42%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0460"
$ws.Range("E45").Value = "  +5.17%  "

$ws.Range("B46").Value = "ThetaToken"
$ws.Range("C46").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.96"
$ws.Range("E46").Value = "  +0.64%  "

$ws.Range("B47").Value = "ApeXProtocol"
$ws.Range("C47").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.45"
$ws.Range("E47").Value = "  -2.19%  "

$ws.Range("E48").Value = "  +3.74%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "9.21"
$ws.Range("E49").Value = "  +2.71%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.00"
$ws.Range("E50").Value = "  +0.05%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.45"
$ws.Range("E51").Value = "  +8.16%  "
